# "New research info. Moved teaching info around."
#
# The sheet lists research/teaching experience entries, one per "block" of
# rows (a data row with what/when/with/where, followed by one or more
# "why" continuation rows in column I). A new "Doctoral research" entry is
# inserted at the very top (pushing every existing entry down), and the
# blank spacer rows that used to separate blocks are squeezed out so the
# remaining blocks sit back-to-back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for the new top entry: insert 4 blank rows before row 2.
#    (Excel copies the formatting of the surrounding rows onto the new
#    blank cells automatically, same as a manual "Insert Rows" in the UI.)
$ws.Rows("2:5").Insert()

# 2) The existing entries used to be separated by a blank spacer row. After
#    the insert above those spacer rows now live at row 13 and row 16 (the
#    second shifts up once the first is removed) - delete them so the
#    blocks are contiguous again.
$ws.Rows("13").Delete()
$ws.Rows("16").Delete()

# The two continuation rows (3 and 4) of the new entry only need column I -
# column H carried a leftover blank/styled cell from the insert because the
# whole column has a style; clear it so those rows stay truly blank there.
$ws.Range("H3:H5").Clear()

# 3) Column A needed a touch more width once the new "Doctoral research"
#    label was added.
$ws.Columns("A").ColumnWidth = 36

# 4) Fill in the new "Doctoral research" entry in the rows just opened up.
$ws.Range("A2").Value = "Doctoral research"
$ws.Range("B2").Value = "August"
$ws.Range("C2").Value = 2018
$ws.Range("D2").Value = "Present"
$ws.Range("G2").Value = "Ecology and Evolutionary Biology"
$ws.Range("H2").Value = "Tulane University, New Orleans, LA"
$ws.Range("I2").Value = "“Leaf Functional Traits Influence on Foliar Endophytic Fungi and Their Effects on Plant’s Response to Herbivory and Pathogenicity In Tropical Trees”  (in progress)"
$ws.Range("I3").Value = "“Foliar Endophytic Fungi in Yellow Monkeyflowers Along an Elevational Gradient in the Sierra Nevada, CA”"
$ws.Range("I4").Value = "“Leaf Functional Trait Plasticity and Foliar Endophytic Fungi in Yellow Monkeyflowers: Linking Traits and Symbionts to Genes”"
$ws.Range("I5").Value = "PI, Sunshine Van Bael Ph.D. & Kathleen Ferris, Ph.D."

# 5) Leave the cursor where the author left it.
$ws.Range("A16").Select() | Out-Null
